# Apply updated ranking data to the active worksheet.
# Rows 6-10 get some of their identity columns (prolificid / name / race)
# swapped around (as the underlying ranking order was recomputed), and the
# realeffort score column (F) is refreshed for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 13.04325402792447

# Row 3
$ws.Range("F3").Value = 8.342377812971202

# Row 4
$ws.Range("F4").Value = 7.489472321657063

# Row 5
$ws.Range("F5").Value = 7.219432926815826

# Row 6
$ws.Range("B6").Value = 22
$ws.Range("C6").Value = "60db4fde6193c50664c9c478"
$ws.Range("D6").Value = "Edosagbe"
$ws.Range("F6").Value = 5.404956080902719
$ws.Range("G6").Value = "Black or African American"

# Row 7
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("D7").Value = "Jamarii"
$ws.Range("F7").Value = 5.203546488046102

# Row 8
$ws.Range("B8").Value = 26
$ws.Range("C8").Value = "5dd671942b033b5ec8bc97b4"
$ws.Range("D8").Value = "Juan"
$ws.Range("F8").Value = 5.194694186643499
$ws.Range("G8").Value = "Hispanic"

# Row 9
$ws.Range("B9").Value = 33
$ws.Range("C9").Value = "60b322994d0b901954690036"
$ws.Range("D9").Value = "Brennan"
$ws.Range("F9").Value = 4.334666484926464

# Row 10
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "5e2522d6b734b47915f88275"
$ws.Range("D10").Value = "Corey"
$ws.Range("F10").Value = 4.178693876440433

# Row 11
$ws.Range("F11").Value = 3.419194189605884

# Row 12
$ws.Range("F12").Value = 2.385885516067507

# Row 13
$ws.Range("F13").Value = 2.223286854337817
